function Set-TextCell {
    param($sheet, $addr, $val)
    $c = $sheet.Range($addr)
    $c.Formula = "'" + $val
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
Set-TextCell $ws "D2" "43.680.83"
$ws.Range("E2").Value = "  -0.04%  "

# Row 3
Set-TextCell $ws "D3" "2.279.94"
$ws.Range("E3").Value = "  -0.41%  "

# Row 4
$ws.Range("E4").Value = "  +0.37%  "

# Row 5
Set-TextCell $ws "D5" "114.01"
$ws.Range("E5").Value = "  +10.12%  "

# Row 6
Set-TextCell $ws "D6" "267.58"
$ws.Range("E6").Value = "  -1.16%  "

# Row 7
Set-TextCell $ws "D7" "0.628"
$ws.Range("E7").Value = "  +0.53%  "

# Row 8
$ws.Range("E8").Value = "  +0.25%  "

# Row 9
Set-TextCell $ws "D9" "0.617"
$ws.Range("E9").Value = "  +1.70%  "

# Row 10
Set-TextCell $ws "D10" "48.86"
$ws.Range("E10").Value = "  +6.51%  "

# Row 11
Set-TextCell $ws "D11" "0.0940"
$ws.Range("E11").Value = "  +0.79%  "

# Row 12
Set-TextCell $ws "D12" "8.93"
$ws.Range("E12").Value = "  +11.95%  "

# Row 13
$ws.Range("E13").Value = "  +0.61%  "

# Row 14
Set-TextCell $ws "D14" "15.69"
$ws.Range("E14").Value = "  +0.70%  "

# Row 15
Set-TextCell $ws "D15" "2.626.85"
$ws.Range("E15").Value = "  -0.24%  "

# Row 16
Set-TextCell $ws "D16" "0.878"
$ws.Range("E16").Value = "  +2.43%  "

# Row 17
Set-TextCell $ws "D17" "2.286.27"
$ws.Range("E17").Value = "  -0.07%  "

# Row 18
Set-TextCell $ws "D18" "43.551.96"
$ws.Range("E18").Value = "  -0.28%  "

# Row 19
$ws.Range("E19").Value = "  -0.79%  "

# Row 20
Set-TextCell $ws "D20" "6.99"
$ws.Range("E20").Value = "  +11.69%  "

# Row 21
Set-TextCell $ws "D21" "72.20"
$ws.Range("E21").Value = "  -0.03%  "

# Row 22
Set-TextCell $ws "D22" "2.40"
$ws.Range("E22").Value = "  -3.22%  "

# Row 23
Set-TextCell $ws "D23" "9.96"
$ws.Range("E23").Value = "  +9.02%  "

# Row 24
Set-TextCell $ws "D24" "232.82"
$ws.Range("E24").Value = "  +0.43%  "

# Row 25
$ws.Range("E25").Value = "  -0.09%  "

# Row 27
Set-TextCell $ws "D27" "11.60"
$ws.Range("E27").Value = "  +3.76%  "

# Row 28
Set-TextCell $ws "D28" "41.94"
$ws.Range("E28").Value = "  +3.83%  "

# Row 29
$ws.Range("E29").Value = "  -1.64%  "

# Row 30
Set-TextCell $ws "D30" "2.25"
$ws.Range("E30").Value = "  +1.05%  "

# Row 31
Set-TextCell $ws "D31" "173.25"
$ws.Range("E31").Value = "  -2.18%  "

# Row 32
Set-TextCell $ws "D32" "21.54"
$ws.Range("E32").Value = "  -1.09%  "

# Row 33
Set-TextCell $ws "D33" "0.0920"
$ws.Range("E33").Value = "  +2.27%  "

# Row 34
Set-TextCell $ws "D34" "5.67"
$ws.Range("E34").Value = "  +3.79%  "

# Row 35
Set-TextCell $ws "D35" "0.127"
$ws.Range("E35").Value = "  +0.26%  "

# Row 36
Set-TextCell $ws "D36" "4.66"
$ws.Range("E36").Value = "  -5.13%  "

# Row 37
Set-TextCell $ws "D37" "0.0354"
$ws.Range("E37").Value = "  +0.31%  "

# Row 38
$ws.Range("E38").Value = "  -3.19%  "

# Row 39
$ws.Range("E39").Value = "  +4.85%  "

# Row 40
Set-TextCell $ws "D40" "14.89"
$ws.Range("E40").Value = "  +20.85%  "

# Row 41
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell $ws "D41" "2.42"
$ws.Range("E41").Value = "  +4.30%  "

# Row 42
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextCell $ws "D42" "74.06"
$ws.Range("E42").Value = "  +12.90%  "

# Row 43
Set-TextCell $ws "D43" "0.238"
$ws.Range("E43").Value = "  +0.24%  "

# Row 44
Set-TextCell $ws "D44" "6.31"
$ws.Range("E44").Value = "  +20.19%  "

# Row 45
$ws.Range("E45").Value = "  +0.23%  "

# Row 46
Set-TextCell $ws "D46" "1.38"
$ws.Range("E46").Value = "  +1.22%  "

# Row 47
Set-TextCell $ws "D47" "8.68"
$ws.Range("E47").Value = "  -1.06%  "

# Row 48
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell $ws "D48" "0.100"
$ws.Range("E48").Value = "  -1.57%  "

# Row 49
Set-TextCell $ws "D49" "102.20"
$ws.Range("E49").Value = "  +3.21%  "

# Row 50
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell $ws "D50" "1.25"
$ws.Range("E50").Value = "  +1.98%  "

# Row 51
Set-TextCell $ws "D51" "0.453"
$ws.Range("E51").Value = "  +1.36%  "
